# Product Backlog and Sprints - apply commit changes:
#  - Rename "Convert Scraping and Analysis Scripts from Jupyter to .py" (A5)
#    to "Convert Scraping Scripts from Jupyter to .py", and mark it Done (F5).
#  - Mark "Refactor analysis scripts..." (A7) Done (F7).
#  - Insert a new task "Convert Thunder analysis scripts from Jupyter to .py"
#    (Sprint 2, Done) right after it.
#  - Move "Update thunder data analysis..." into Sprint 3 (now right below the
#    new Thunder-scripts task).
#  - Insert a new placeholder task "Convert Meteorite analysis scripts from
#    Jupyter to .py" right after that.
#  - Add a "Comments" column header (G1).
#  - Update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the scraping-scripts task and mark it + row 7 as Done ---
$ws.Range("A5").Value = "Convert Scraping Scripts from Jupyter to .py"

$ws.Range("F2").Copy()
$ws.Range("F5").PasteSpecial(-4122)
$ws.Range("F5").Value = 45918

$ws.Range("F2").Copy()
$ws.Range("F7").PasteSpecial(-4122)
$ws.Range("F7").Value = 45918

# --- Make room for two new rows: one right after row 7 (Thunder scripts),
#     one after the (soon to be relocated) thunder-data-analysis row. ---
$ws.Rows.Item(8).Resize(3).Insert()

# Row 8: Convert Thunder analysis scripts from Jupyter to .py
$ws.Range("A8").Value = "Convert Thunder analysis scripts from Jupyter to .py"
$ws.Range("B8").Value = "$$$"
$ws.Range("C8").Value = 13
$ws.Range("D8").Value = "No"
$ws.Range("E8").Value = 2
$ws.Range("F2").Copy()
$ws.Range("F8").PasteSpecial(-4122)
$ws.Range("F8").Value = 45918

# Row 9: relocated "Update thunder data analysis..." task, now Sprint 3
$ws.Range("A9").Value = "Update thunder data analysis to use a centered moving 10 year average"
$ws.Range("B9").Value = "$$$"
$ws.Range("C9").Value = 21
$ws.Range("D9").Value = "No"
$ws.Range("E9").Value = 3

# Row 10: new placeholder task for the meteorite analysis conversion
$ws.Range("A10").Value = "Convert Meteorite analysis scripts from Jupyter to .py"
$ws.Range("B10").Value = "$$$"
$ws.Range("C10").Value = 21
$ws.Range("D10").Value = "No"

# --- Remove the old copy of the thunder-data-analysis row, now at row 17 ---
$ws.Rows.Item(17).Delete()

# --- New "Comments" column header ---
$ws.Range("B1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "Comments"
$ws.Columns.Item(7).ColumnWidth = 9.81640625

# --- Selection, to match the saved UI state ---
$ws.Range("E10").Select()
